$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 1621.375
$ws.Cells.Item(80, 9).Value = 2167.25
$ws.Cells.Item(80, 11).Value = 6501.75
$ws.Cells.Item(80, 13).Value = -5503.75
$ws.Cells.Item(83, 8).Value = 1621.375
$ws.Cells.Item(83, 9).Value = 2167.25
$ws.Cells.Item(83, 11).Value = 19505.25
$ws.Cells.Item(83, 13).Value = -14513.25
$ws.Cells.Item(88, 8).Value = 3624.8125
$ws.Cells.Item(88, 9).Value = 1200
$ws.Cells.Item(88, 10).Value = 4433.0835
$ws.Cells.Item(88, 11).Value = 1200
$ws.Cells.Item(88, 12).Value = 4433.0835
$ws.Cells.Item(88, 13).Value = -794
$ws.Cells.Item(88, 14).Value = -5245.0835
$ws.Cells.Item(91, 8).Value = 3624.8125
$ws.Cells.Item(91, 9).Value = 1200
$ws.Cells.Item(91, 10).Value = 4433.0835
$ws.Cells.Item(91, 11).Value = 1200
$ws.Cells.Item(91, 12).Value = 4433.0835
$ws.Cells.Item(91, 13).Value = 204
$ws.Cells.Item(91, 14).Value = -7241.0835
$ws.Cells.Item(98, 9).Value = 6472
$ws.Cells.Item(98, 10).Value = 5149.5
$ws.Cells.Item(98, 11).Value = 6472
$ws.Cells.Item(98, 12).Value = 5149.5
$ws.Cells.Item(98, 13).Value = -4974
$ws.Cells.Item(98, 14).Value = -8145.5
$ws.Cells.Item(112, 8).Value = 2046.5186
$ws.Cells.Item(112, 10).Value = 2082.923
$ws.Cells.Item(112, 12).Value = 6248.768999999999
$ws.Cells.Item(112, 14).Value = -8464.769
$ws.Cells.Item(122, 9).Value = 6472
$ws.Cells.Item(122, 10).Value = 5149.5
$ws.Cells.Item(122, 11).Value = 19416
$ws.Cells.Item(122, 12).Value = 15448.5
$ws.Cells.Item(122, 13).Value = -16966
$ws.Cells.Item(122, 14).Value = -20348.5
$ws.Cells.Item(127, 8).Value = 2446.25
$ws.Cells.Item(127, 9).Value = 2528.3333
$ws.Cells.Item(127, 11).Value = 7584.999899999999
$ws.Cells.Item(127, 13).Value = -2624.999899999999
$ws.Cells.Item(131, 8).Value = 2301.7896
$ws.Cells.Item(131, 9).Value = 773.1667
$ws.Cells.Item(131, 10).Value = 4922.2856
$ws.Cells.Item(131, 11).Value = 2319.5001
$ws.Cells.Item(131, 12).Value = 14766.8568
$ws.Cells.Item(131, 13).Value = 2720.4999
$ws.Cells.Item(131, 14).Value = -24846.8568
$ws.Cells.Item(132, 8).Value = 6411586.5
$ws.Cells.Item(132, 9).Value = 7247608.5
$ws.Cells.Item(132, 11).Value = 21742825.5
$ws.Cells.Item(132, 13).Value = -21740295.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2639.6204
$ws.Cells.Item(32, 9).Value = 2322.8289
$ws.Cells.Item(32, 11).Value = 2322.8289
$ws.Cells.Item(32, 13).Value = -2035.8289
$ws.Cells.Item(61, 8).Value = 71432180
$ws.Cells.Item(61, 9).Value = 55558224
$ws.Cells.Item(61, 11).Value = 55558224
$ws.Cells.Item(61, 13).Value = -55558012
$ws.Cells.Item(74, 8).Value = 2177.7058
$ws.Cells.Item(74, 9).Value = 1826.6666
$ws.Cells.Item(74, 10).Value = 2572.625
$ws.Cells.Item(74, 11).Value = 1826.6666
$ws.Cells.Item(74, 12).Value = 2572.625
$ws.Cells.Item(74, 13).Value = -952.6666
$ws.Cells.Item(74, 14).Value = -4320.625
$ws.Cells.Item(77, 8).Value = 2177.7058
$ws.Cells.Item(77, 9).Value = 1826.6666
$ws.Cells.Item(77, 10).Value = 2572.625
$ws.Cells.Item(77, 11).Value = 9133.333000000001
$ws.Cells.Item(77, 12).Value = 12863.125
$ws.Cells.Item(77, 13).Value = -4765.333000000001
$ws.Cells.Item(77, 14).Value = -21599.125
$ws.Cells.Item(102, 8).Value = 1081
$ws.Cells.Item(102, 9).Value = 1081
$ws.Cells.Item(102, 11).Value = 1081
$ws.Cells.Item(102, 13).Value = 541
$ws.Cells.Item(132, 8).Value = 1446.9592
$ws.Cells.Item(132, 9).Value = 1091.2903
$ws.Cells.Item(132, 11).Value = 3273.8709
$ws.Cells.Item(132, 13).Value = -743.8708999999999
$ws.Cells.Item(136, 8).Value = 71432180
$ws.Cells.Item(136, 9).Value = 55558224
$ws.Cells.Item(136, 11).Value = 166674672
$ws.Cells.Item(136, 13).Value = -166672122

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 4083.611
$ws.Cells.Item(122, 10).Value = 4981.222
$ws.Cells.Item(122, 12).Value = 14943.666
$ws.Cells.Item(122, 14).Value = -19843.666
$ws.Cells.Item(132, 8).Value = 1588.5682
$ws.Cells.Item(132, 9).Value = 1164.6052
$ws.Cells.Item(132, 11).Value = 3493.8156
$ws.Cells.Item(132, 13).Value = -963.8155999999999
$ws.Cells.Item(134, 8).Value = 1979.6316
$ws.Cells.Item(134, 9).Value = 1775.4062
$ws.Cells.Item(134, 10).Value = 3068.8333
$ws.Cells.Item(134, 11).Value = 5326.2186
$ws.Cells.Item(134, 12).Value = 9206.499899999999
$ws.Cells.Item(134, 13).Value = -2791.2186
$ws.Cells.Item(134, 14).Value = -14276.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 5962922.5
$ws.Cells.Item(131, 9).Value = 71429130
$ws.Cells.Item(131, 10).Value = 11449.767
$ws.Cells.Item(131, 11).Value = 214287390
$ws.Cells.Item(131, 12).Value = 34349.301
$ws.Cells.Item(131, 13).Value = -214282350
$ws.Cells.Item(131, 14).Value = -44429.301

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1747.8182
$ws.Cells.Item(122, 9).Value = 1672.4286
$ws.Cells.Item(122, 11).Value = 5017.2858
$ws.Cells.Item(122, 13).Value = -2567.2858
$ws.Cells.Item(132, 8).Value = 1925863.2
$ws.Cells.Item(132, 9).Value = 3498770.2
$ws.Cells.Item(132, 11).Value = 10496310.6
$ws.Cells.Item(132, 13).Value = -10493780.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3687.8333
$ws.Cells.Item(7, 9).Value = 3170.8
$ws.Cells.Item(7, 11).Value = 3170.8
$ws.Cells.Item(7, 13).Value = -3058.8
$ws.Cells.Item(22, 8).Value = 2592.2222
$ws.Cells.Item(22, 9).Value = 5300
$ws.Cells.Item(22, 10).Value = 1818.5714
$ws.Cells.Item(22, 11).Value = 5300
$ws.Cells.Item(22, 12).Value = 1818.5714
$ws.Cells.Item(22, 13).Value = -5005
$ws.Cells.Item(22, 14).Value = -2408.5714
$ws.Cells.Item(27, 8).Value = 2592.2222
$ws.Cells.Item(27, 9).Value = 5300
$ws.Cells.Item(27, 10).Value = 1818.5714
$ws.Cells.Item(27, 11).Value = 5300
$ws.Cells.Item(27, 12).Value = 1818.5714
$ws.Cells.Item(27, 13).Value = -5193
$ws.Cells.Item(27, 14).Value = -2032.5714
$ws.Cells.Item(82, 8).Value = 2582.1428
$ws.Cells.Item(82, 9).Value = 1734
$ws.Cells.Item(82, 10).Value = 3218.25
$ws.Cells.Item(82, 11).Value = 1734
$ws.Cells.Item(82, 12).Value = 3218.25
$ws.Cells.Item(82, 13).Value = -1373
$ws.Cells.Item(82, 14).Value = -3940.25
$ws.Cells.Item(85, 8).Value = 2582.1428
$ws.Cells.Item(85, 9).Value = 1734
$ws.Cells.Item(85, 10).Value = 3218.25
$ws.Cells.Item(85, 11).Value = 1734
$ws.Cells.Item(85, 12).Value = 3218.25
$ws.Cells.Item(85, 13).Value = -486
$ws.Cells.Item(85, 14).Value = -5714.25
$ws.Cells.Item(126, 8).Value = 3687.8333
$ws.Cells.Item(126, 9).Value = 3170.8
$ws.Cells.Item(126, 11).Value = 9512.400000000001
$ws.Cells.Item(126, 13).Value = -7042.400000000001
$ws.Cells.Item(136, 8).Value = 3111.4827
$ws.Cells.Item(136, 9).Value = 2024.2174
$ws.Cells.Item(136, 10).Value = 7279.3335
$ws.Cells.Item(136, 11).Value = 6072.6522
$ws.Cells.Item(136, 12).Value = 21838.0005
$ws.Cells.Item(136, 13).Value = -3522.6522
$ws.Cells.Item(136, 14).Value = -26938.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 73248.55
$ws.Cells.Item(122, 9).Value = 89026
$ws.Cells.Item(122, 11).Value = 267078
$ws.Cells.Item(122, 13).Value = -264628
$ws.Cells.Item(132, 8).Value = 1471.3948
$ws.Cells.Item(132, 9).Value = 1055.4286
$ws.Cells.Item(132, 11).Value = 3166.2858
$ws.Cells.Item(132, 13).Value = -636.2857999999997
$ws.Cells.Item(136, 8).Value = 10685654
$ws.Cells.Item(136, 9).Value = 13552068
$ws.Cells.Item(136, 10).Value = 1746.4546
$ws.Cells.Item(136, 11).Value = 40656204
$ws.Cells.Item(136, 12).Value = 5239.3638
$ws.Cells.Item(136, 13).Value = -40653654
$ws.Cells.Item(136, 14).Value = -10339.3638
